# The commit swaps the content of ppt/theme/theme1.xml ("Office Theme")
# and ppt/theme/theme2.xml ("Integral") so that the deck's active theme
# (the one backing the slide master / slides, physically stored at
# ppt/theme/theme2.xml) ends up carrying the plain "Office" colour
# scheme instead of the custom "Integral" one (and vice versa for the
# file that used to hold the Office colours).
#
# The PowerPoint object model only exposes the *values* of a theme's
# colour scheme for editing (ThemeColorScheme(i).RGB) -- there is no
# settable "theme name" / "rename part" verb on Theme / Design /
# ColorScheme, and every master/notes-master/handout-master/slide in
# this deck resolves to the same single active Theme object (the one
# stored in ppt/theme/theme2.xml). So we reproduce the swap's visible
# effect by rewriting that theme's 12 colour-scheme entries, in OOXML
# clrScheme order (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink), to
# the "Office Theme" values that theme1.xml held before the edit.

function HexToOleColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$colorScheme = $master.Theme.ThemeColorScheme

# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Item($i).RGB = HexToOleColor $officeThemeColors[$i - 1]
}
